# Auto-generated: update crypto price/volume figures per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value() = "'26.172.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value() = "'  +3.79%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value() = "'1.604.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value() = "'  +3.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value() = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value() = "'212.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value() = "'  +2.86%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value() = "'  -0.25%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value() = "'  +2.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value() = "'  +3.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value() = "'  +1.63%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value() = "'18.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value() = "'  +1.47%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value() = "'0.0819"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value() = "'  +5.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value() = "'1.826.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value() = "'  +3.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value() = "'1.601.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value() = "'  +3.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value() = "'  +0.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value() = "'  +1.81%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value() = "'26.154.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value() = "'  +3.82%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value() = "'60.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value() = "'  +3.32%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value() = "'0.0₃0723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value() = "'  +2.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value() = "'205.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value() = "'  +10.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value() = "'  +3.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value() = "'  +0.93%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value() = "'  +2.86%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value() = "'1.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value() = "'  +10.79%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value() = "'141.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value() = "'  +1.50%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value() = "'  -0.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value() = "'  -2.87%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value() = "'  +2.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value() = "'  +0.94%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value() = "'  +1.91%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value() = "'  +2.20%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value() = "'  +3.74%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value() = "'  +0.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value() = "'1.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value() = "'  +1.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value() = "'  +1.98%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value() = "'  +10.46%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value() = "'1.117.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value() = "'  +3.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value() = "'  -0.23%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value() = "'ARBITRUM"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value() = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value() = "'0.783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value() = "'  +3.43%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value() = "'MXToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value() = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value() = "'2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value() = "'  +1.85%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value() = "'0.492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value() = "'  -0.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value() = "'  -2.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value() = "'1.738.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value() = "'  +3.43%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value() = "'  +1.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value() = "'92.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value() = "'  +0.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value() = "'1.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value() = "'  +5.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value() = "'53.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value() = "'  +2.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E49").Value() = "'  +1.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value() = "'  -0.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value() = "'0.0₇0926"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value() = "'  -17.02%  "
$ws.Range("E51").Style = "Normal"
